$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old rows 14-19 (data no longer present after the edit)
$ws.Rows("14:19").Delete()

# New data for rows 2-13 (id, title, price, date, city, description)
$data = @(
    @("10132208", "Сдаются в аренду места для бровиста/визажиста", 11000, "Сегодня, 15:31", "Севастополь", "Тип объявления: Сдам; Вид объекта: Салон красоты; Готовый бизнес: нет;"),
    @("9990054", "1-к, Дмитрия Ульянова-2, Стрелецкая бухта, Гагаринский район.", 25000, "Сегодня, 13:35", "Севастополь", "Тип объявления: Сдам; Срок аренды: На длительный срок; Количество комнат: 1; Тип дома: Блочный; Площадь: 36; Этаж: 3; Этажей в доме: 5;"),
    @("8652548", "1-к квартира, Дмитрия Ульянова-2, 36 м2, 3/5 эт.", 25000, "Сегодня, 13:35", "Севастополь", "Тип объявления: Сдам; Срок аренды: На длительный срок; Количество комнат: 1; Тип дома: Блочный; Площадь: 36; Этаж: 3; Этажей в доме: 5;"),
    @("10132004", "Арендуй выгодно офис в Евпатории", 10000, "Сегодня, 13:25", "Евпатория", "Тип объявления: Сдам; Вид объекта: Офис; Готовый бизнес: нет;"),
    @("5307964", "Сдам нежилое помещение", 25000, "Сегодня, 08:47", "Симферополь", "Тип объявления: Сдам; Вид объекта: Помещение; Готовый бизнес: нет;"),
    @("7105101", "Сдам нежилое помещение.", 23000, "Сегодня, 08:47", "Симферополь", "Тип объявления: Сдам; Вид объекта: Помещение; Готовый бизнес: да;"),
    @("1864998", "Аренда магазина в Крыму на берегу Черного моря.", 35000, "Сегодня, 08:01", "Симферополь", "Тип объявления: Сдам; Вид объекта: Магазин; Готовый бизнес: да;"),
    @("10096859", "Аренда 2 х комнатной квартиры Симферопольский район село Мирное", 30000, "Сегодня, 06:09", "Симферополь", "Тип объявления: Сдам; Срок аренды: На длительный срок; Количество комнат: 2; Площадь: 46; Этаж: 1; Этажей в доме: 1;"),
    @("10052833", "Аренда складского помещение ул Индустриальная /ул Лексина", 36000, "Сегодня, 06:08", "Симферополь", "Тип объявления: Сдам; Вид объекта: Склад;"),
    @("9677392", "Аренда помещение свободного назначение", 30000, "Сегодня, 06:06", "Симферополь", "Тип объявления: Сдам; Вид объекта: Салон красоты;"),
    @("9150844", "Аренда торгового помещение", 45000, "Сегодня, 06:05", "Симферополь", "Тип объявления: Сдам; Вид объекта: Магазин; Готовый бизнес: да;"),
    @("5847411", "Аренда помещение свободного назначение", 40000, "Сегодня, 06:03", "Симферополь", "Тип объявления: Сдам; Вид объекта: Помещение; Готовый бизнес: да;")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
